$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. Delete the empty spacer paragraph right before the "TimeLess" bullet
#    (w:spacing line=220 exact, w:ind left=720, empty run props sz=18).
# ---------------------------------------------------------------------------
$d.Paragraphs(9).Range.Delete()

# After the deletion every following paragraph index shifts down by one:
#   old #12 (Relevant Courses)        -> #11
#   old #13 (empty -> Certifications) -> #12
#   old #49 (little overhead...)      -> #48
#   old #52 (OAG Analytics / _GoBack) -> #51
#   old #62 (hobbies)                 -> #61
#   old #63 (black belt)              -> #62
#   old #64 (National Honor Society)  -> #63

# ---------------------------------------------------------------------------
# 2. "Relevant Courses" paragraph: split the long run, change "Development"
#    to "Dev.", and append ", Big Data" at the end.
# ---------------------------------------------------------------------------
$coursesXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Relevant Courses:</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">Data Structures, Data Systems, </w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Operating Systems</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>, Database Systems, Algorithms, Software</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> Dev</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>.</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>, Data Mining</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>, Big Data</w:t></w:r>' + `
  '</w:p>'
$d.Paragraphs(11).Range.InsertXML($coursesXml)

# ---------------------------------------------------------------------------
# 3. Remove the old _GoBack bookmark that used to sit before "June 2015"
#    *before* re-adding it elsewhere below - "_GoBack" is a singleton
#    bookmark, so the stale copy has to go first or the name lookup keeps
#    resolving to whichever one is newest.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 4. Previously-empty paragraph after "Relevant Courses": add the new
#    "Relevant Certifications" line, finishing with a _GoBack bookmark.
# ---------------------------------------------------------------------------
$certsXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:i/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Relevant Certifications:</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>IBM Cognitive Class</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>es</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>:</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> Data Science 101, Data Science Methodology, Python for Data Science,</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$d.Paragraphs(12).Range.InsertXML($certsXml)

# ---------------------------------------------------------------------------
# 5. Drop the trailing period after "with little overhead supervision." -
#    rebuild the paragraph keeping its original two-run split (editing the
#    run text in place rather than Find/Replace, which would merge the runs).
# ---------------------------------------------------------------------------
$overheadXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="BodyText"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Reached out to homeowners, visited properties, and completed thorough property walks</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> with little overhead supervision</w:t></w:r>' + `
  '</w:p>'
$d.Paragraphs(48).Range.InsertXML($overheadXml)

# ---------------------------------------------------------------------------
# 6-8. Rotate the three "ADDITIONAL INFORMATION" bullet paragraphs: a brand
#    new "Data Science Team / Cryptocurrency Club" bullet is inserted first,
#    pushing the hobbies + black-belt bullets down one slot each, and the
#    old "National Honor Society" bullet is dropped.
# ---------------------------------------------------------------------------
$clubsXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="BodyText"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">Part of the CU Boulder Data Science </w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Team</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> and Cryptocurrency Club </w:t></w:r>' + `
  '</w:p>'

$hobbiesXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="BodyText"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">I enjoy playing </w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>intermural football,</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">hiking, </w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>wakeboarding,</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> snowmobilin</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>g, and actively participating</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> in both wrestling and rugby</w:t></w:r>' + `
  '</w:p>'

$blackBeltXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="BodyText"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>R</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>ecei</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>ved a black belt in Tae Kwon Do</w:t></w:r>' + `
  '</w:p>'

$d.Paragraphs(61).Range.InsertXML($clubsXml)
$d.Paragraphs(62).Range.InsertXML($hobbiesXml)
$d.Paragraphs(63).Range.InsertXML($blackBeltXml)
